# Rename the three header/footer logo pictures:
#   - Pearson logo in the "first page" footer  (footer1.xml, docPr id="3")
#   - Pearson logo in the "default" footer     (footer2.xml, docPr id="2")
#   - BTEC logo   in the "first page" header   (header1.xml, docPr id="1")
# went from image1.png/image1.png/image2.jpg -> image2.png/image2.png/image1.jpg
# (the BTEC logo and the Pearson logos effectively swapped their generic
# "imageN" file-name labels).

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footers (Pearson logo, 952500 x 285750 EMU) -------------------------
# Footers.Item(1) = wdHeaderFooterPrimary -> footer2.xml (docPr id="2")
$footerPrimary = $sec.Footers.Item(1)
$footerPrimary.Range.InlineShapes.Item(1).Name = "image2.png"

# Footers.Item(2) = wdHeaderFooterFirstPage -> footer1.xml (docPr id="3")
$footerFirstPage = $sec.Footers.Item(2)
$footerFirstPage.Range.InlineShapes.Item(1).Name = "image2.png"

# --- Header (BTEC logo, 914400 x 277792 EMU) ------------------------------
# Headers.Item(2) = wdHeaderFooterFirstPage -> header1.xml (docPr id="1")
$headerFirstPage = $sec.Headers.Item(2)
$headerFirstPage.Range.InlineShapes.Item(1).Name = "image1.jpg"
